$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-01-28 04:57:34"
$wsZh.Range("G2").Value = "2016-01-28 04:58:14"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-01-28 04:57:44"
$wsDe.Range("G2").Value = "2016-01-28 04:58:31"
